{"js": "// Version-control table, first data row (version \"2.8.1\"): swap the\n// \"responsible\" (\u0e1c\u0e39\u0e49\u0e23\u0e31\u0e1a\u0e1c\u0e34\u0e14\u0e0a\u0e2d\u0e1a) and \"reviewer\" (\u0e1c\u0e39\u0e49\u0e15\u0e23\u0e27\u0e08) names.\n//   \u0e1c\u0e39\u0e49\u0e23\u0e31\u0e1a\u0e1c\u0e34\u0e14\u0e0a\u0e2d\u0e1a: \"\u0e27\u0e34\u0e23\u0e31\u0e15\u0e19\u0e4c (TL)\"    -> \"\u0e13\u0e31\u0e10\u0e14\u0e19\u0e31\u0e22 (DM)\"\n//   \u0e1c\u0e39\u0e49\u0e15\u0e23\u0e27\u0e08:      \"\u0e01\u0e34\u0e15\u0e15\u0e34\u0e1e\u0e28 (SP)\"   -> \"\u0e27\u0e34\u0e23\u0e31\u0e15\u0e19\u0e4c (TL)\"\nconst table = context.document.body.tables.getFirst();\n\n// Row index 5 (0-based): \u0e0a\u0e37\u0e48\u0e2d\u0e40\u0e2d\u0e01\u0e2a\u0e32\u0e23(0), \u0e40\u0e27\u0e2d\u0e23\u0e4c\u0e0a\u0e31\u0e19\u0e1b\u0e31\u0e08\u0e08\u0e38\u0e1a\u0e31\u0e19(1), \u0e27\u0e31\u0e19\u0e17\u0e35\u0e48\u0e41\u0e01\u0e49\u0e44\u0e02\u0e25\u0e48\u0e32\u0e2a\u0e38\u0e14(2),\n// \u0e1c\u0e39\u0e49\u0e04\u0e27\u0e1a\u0e04\u0e38\u0e21\u0e40\u0e27\u0e2d\u0e23\u0e4c\u0e0a\u0e31\u0e19(3), header row(4), first version row \"2.8.1\"(5).\n// Column index 3 = \u0e1c\u0e39\u0e49\u0e23\u0e31\u0e1a\u0e1c\u0e34\u0e14\u0e0a\u0e2d\u0e1a, column index 4 = \u0e1c\u0e39\u0e49\u0e15\u0e23\u0e27\u0e08.\nconst respCell = table.getCell(5, 3);\nconst respPara = respCell.body.paragraphs.getFirst();\n\nconst reviewCell = table.getCell(5, 4);\nconst reviewPara = reviewCell.body.paragraphs.getFirst();\n\n// --- \u0e1c\u0e39\u0e49\u0e23\u0e31\u0e1a\u0e1c\u0e34\u0e14\u0e0a\u0e2d\u0e1a cell: \"\u0e27\u0e34\u0e23\u0e31\u0e15\u0e19\u0e4c\" + \" (TL)\" -> \"\u0e13\u0e31\u0e10\u0e14\u0e19\u0e31\u0e22\" + \" (DM)\" ---\nconst respName = respPara.search(\"\u0e27\u0e34\u0e23\u0e31\u0e15\u0e19\u0e4c\", { matchCase: true });\nrespName.load(\"items\");\nawait context.sync();\nrespName.items[0].insertText(\"\u0e13\u0e31\u0e10\u0e14\u0e19\u0e31\u0e22\", Word.InsertLocation.replace);\n\nconst respRole = respPara.search(\"(TL)\", { matchCase: true });\nrespRole.load(\"items\");\nawait context.sync();\nrespRole.items[0].insertText(\"(DM)\", Word.InsertLocation.replace);\nawait context.sync();\n\n// --- \u0e1c\u0e39\u0e49\u0e15\u0e23\u0e27\u0e08 cell: \"\u0e01\u0e34\u0e15\u0e15\u0e34\u0e1e\u0e28 \" + \"(SP)\" -> \"\u0e27\u0e34\u0e23\u0e31\u0e15\u0e19\u0e4c\" + \" (TL)\" ---\nconst reviewName = reviewPara.search(\"\u0e01\u0e34\u0e15\u0e15\u0e34\u0e1e\u0e28 \", { matchCase: true });\nreviewName.load(\"items\");\nawait context.sync();\nreviewName.items[0].insertText(\"\u0e27\u0e34\u0e23\u0e31\u0e15\u0e19\u0e4c\", Word.InsertLocation.replace);\n\nconst reviewRole = reviewPara.search(\"(SP)\", { matchCase: true });\nreviewRole.load(\"items\");\nawait context.sync();\nreviewRole.items[0].insertText(\" (TL)\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Version-control table, first data row (version \"2.8.1\"): swap the\n# \"responsible\" (\u0e1c\u0e39\u0e49\u0e23\u0e31\u0e1a\u0e1c\u0e34\u0e14\u0e0a\u0e2d\u0e1a) and \"reviewer\" (\u0e1c\u0e39\u0e49\u0e15\u0e23\u0e27\u0e08) names.\n#   \u0e1c\u0e39\u0e49\u0e23\u0e31\u0e1a\u0e1c\u0e34\u0e14\u0e0a\u0e2d\u0e1a: \"\u0e27\u0e34\u0e23\u0e31\u0e15\u0e19\u0e4c (TL)\"    -> \"\u0e13\u0e31\u0e10\u0e14\u0e19\u0e31\u0e22 (DM)\"\n#   \u0e1c\u0e39\u0e49\u0e15\u0e23\u0e27\u0e08:      \"\u0e01\u0e34\u0e15\u0e15\u0e34\u0e1e\u0e28 (SP)\"   -> \"\u0e27\u0e34\u0e23\u0e31\u0e15\u0e19\u0e4c (TL)\"\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nfunction Replace-InCellRange($cell, $findText, $replaceText) {\n    $rng = $cell.Range\n    $rng.End = $rng.End - 1   # drop the end-of-cell mark\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $findText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $replaceText\n    $find.Execute([ref]$null, [ref]$true, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, 2) | Out-Null\n}\n\n# Row 6 (1-based): \u0e0a\u0e37\u0e48\u0e2d\u0e40\u0e2d\u0e01\u0e2a\u0e32\u0e23(1), \u0e40\u0e27\u0e2d\u0e23\u0e4c\u0e0a\u0e31\u0e19\u0e1b\u0e31\u0e08\u0e08\u0e38\u0e1a\u0e31\u0e19(2), \u0e27\u0e31\u0e19\u0e17\u0e35\u0e48\u0e41\u0e01\u0e49\u0e44\u0e02\u0e25\u0e48\u0e32\u0e2a\u0e38\u0e14(3),\n# \u0e1c\u0e39\u0e49\u0e04\u0e27\u0e1a\u0e04\u0e38\u0e21\u0e40\u0e27\u0e2d\u0e23\u0e4c\u0e0a\u0e31\u0e19(4), header row(5), first version row \"2.8.1\"(6).\n# Column 4 = \u0e1c\u0e39\u0e49\u0e23\u0e31\u0e1a\u0e1c\u0e34\u0e14\u0e0a\u0e2d\u0e1a, column 5 = \u0e1c\u0e39\u0e49\u0e15\u0e23\u0e27\u0e08.\n$respCell = $t.Cell(6, 4)\nReplace-InCellRange $respCell \"\u0e27\u0e34\u0e23\u0e31\u0e15\u0e19\u0e4c\" \"\u0e13\u0e31\u0e10\u0e14\u0e19\u0e31\u0e22\"\nReplace-InCellRange $respCell \"(TL)\" \"(DM)\"\n\n$reviewCell = $t.Cell(6, 5)\nReplace-InCellRange $reviewCell \"\u0e01\u0e34\u0e15\u0e15\u0e34\u0e1e\u0e28 \" \"\u0e27\u0e34\u0e23\u0e31\u0e15\u0e19\u0e4c\"\nReplace-InCellRange $reviewCell \"(SP)\" \" (TL)\"\n"}
